# Update generated at 456a3b4 (gh-pages): refresh "想去人数"/"最低票价" counters
# across all four sheets, and insert the new "次元格子动漫展" exhibition row
# into Sheet 1 (展览) ahead of the existing "华盟次元动漫嘉年华" / "鸢飞鱼跃" rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet 1 (展览): insert a new row at 40, pushing the two rows below it down ---
$ws1.Rows.Item(40).Insert()

# Copy formatting (bold / border / centered) from the existing numbering column
# down into the freshly inserted A40 cell, matching the style used by A1:A39.
$ws1.Range("A39").Copy()
$ws1.Range("A40").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

# New exhibition occupies row 40 now
$ws1.Range("A40").Value = 39
$ws1.Range("B40").NumberFormat = "@"
$ws1.Range("B40").Value = "2024-10-01"
$ws1.Range("C40").Value = "杭州·第二届次元格子动漫展-进入格子空间，探索次元世界！"
$ws1.Range("D40").Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws1.Range("E40").Value = "2024.10.01 09:30-10.03 17:00"
$ws1.Range("F40").Value = 10
$ws1.Range("G40").Value = 49.9
$ws1.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=90057"
$ws1.Range("I40").Value = "//i2.hdslb.com/bfs/openplatform/202407/Zk5evnyA1722331816981.jpeg"

# Rows 41/42 already hold the shifted-down "华盟" / "鸢飞鱼跃" data (B/C/D/E/H/I untouched);
# only the serial number (A) and the refreshed counters (F/G) need updating.
$ws1.Range("A41").Value = 40
$ws1.Range("F41").Value = 49
$ws1.Range("G41").Value = 60

$ws1.Range("A42").Value = 41
$ws1.Range("F42").Value = 469
$ws1.Range("G42").Value = 85

# --- Sheet 1 numeric refreshes ---
$ws1.Range("F2").Value = 901
$ws1.Range("F3").Value = 834
$ws1.Range("F4").Value = 4412
$ws1.Range("F5").Value = 317
$ws1.Range("F6").Value = 424
$ws1.Range("F7").Value = 3516
$ws1.Range("F8").Value = 998
$ws1.Range("F10").Value = 1346
$ws1.Range("F11").Value = 315
$ws1.Range("G11").Value = 58
$ws1.Range("F12").Value = 311
$ws1.Range("F13").Value = 2383
$ws1.Range("F14").Value = 1256
$ws1.Range("F15").Value = 31
$ws1.Range("F16").Value = 1973
$ws1.Range("F17").Value = 5
$ws1.Range("F18").Value = 521
$ws1.Range("F20").Value = 57
$ws1.Range("F21").Value = 9913
$ws1.Range("F22").Value = 5940
$ws1.Range("F23").Value = 380
$ws1.Range("F26").Value = 133
$ws1.Range("F27").Value = 833
$ws1.Range("F28").Value = 3522
$ws1.Range("F30").Value = 963
$ws1.Range("F31").Value = 456
$ws1.Range("F32").Value = 110
$ws1.Range("F33").Value = 234
$ws1.Range("G33").Value = 158
$ws1.Range("F35").Value = 219
$ws1.Range("F36").Value = 4798
$ws1.Range("F38").Value = 1065
$ws1.Range("F39").Value = 141

# --- Sheet 2 numeric refreshes ---
$ws2.Range("F12").Value = 123
$ws2.Range("F15").Value = 3520
$ws2.Range("F23").Value = 11

# --- Sheet 3 numeric refreshes ---
$ws3.Range("F2").Value = 8683
$ws3.Range("F4").Value = 1549

# --- Sheet 4 numeric refreshes ---
$ws4.Range("F2").Value = 8683
$ws4.Range("F3").Value = 901
$ws4.Range("F4").Value = 1549
$ws4.Range("F6").Value = 4412
$ws4.Range("F8").Value = 424
$ws4.Range("F9").Value = 3516
$ws4.Range("F10").Value = 998
$ws4.Range("F13").Value = 2383
$ws4.Range("F18").Value = 1256
$ws4.Range("F20").Value = 31
$ws4.Range("F21").Value = 123
$ws4.Range("F22").Value = 521
$ws4.Range("F24").Value = 57
$ws4.Range("F25").Value = 9913
$ws4.Range("F26").Value = 3520
$ws4.Range("F28").Value = 380
$ws4.Range("F31").Value = 133
$ws4.Range("F32").Value = 833
$ws4.Range("F33").Value = 3522
$ws4.Range("F35").Value = 963
$ws4.Range("F36").Value = 456
$ws4.Range("F37").Value = 110
$ws4.Range("F40").Value = 219
$ws4.Range("F41").Value = 4798
$ws4.Range("F42").Value = 1065
$ws4.Range("F44").Value = 49
$ws4.Range("F45").Value = 469
$ws4.Range("F47").Value = 11
